# Update operational cost data with the latest P&L calculation.
# Source data: new flight roster (rows 2-26), replacing the previous
# 10-row table and growing it to 25 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    ,@("FL1013", 17742.42, 23754.26, 0,        38853.66, 80350.34)
    ,@("FL1198", 3939.37,  4187.79,  0,        10907.39, 19034.55)
    ,@("FL1294", 8703.97,  9490.78,  0,        18217.4,  36412.15)
    ,@("FL1344", 4029.28,  4793.33,  0,        11585.27, 20407.88)
    ,@("FL1524", 8089.73,  7948.91,  0,        17342.83, 33381.47)
    ,@("FL1541", 36322.68, 77684.29, 0,        48925.23, 162932.2)
    ,@("FL1573", 8951.61,  9490.78,  0,        15621.55, 34063.94)
    ,@("FL1740", 4019.08,  5467.01,  0,        13367.23, 22853.32)
    ,@("FL1832", 17571.6,  28153.2,  3660.75,  34301.88, 83687.43)
    ,@("FL1976", 7231.59,  7591.66,  0,        14062.01, 28885.26)
    ,@("FL1990", 17826.48, 23754.26, 0,        36475.11, 78055.85)
    ,@("FL2240", 8872.98,  13330.89, 0,        16778.86, 38982.73)
    ,@("FL2292", 36715.6,  55306.48, 0,        38329.38, 130351.46)
    ,@("FL2625", 3060.76,  2780.13,  0,        11234.83, 17075.72)
    ,@("FL2644", 39518.1,  59403.25, 0,        38637.52, 137558.87)
    ,@("FL2664", 4057.77,  4793.33,  0,        13174.22, 22025.32)
    ,@("FL2734", 17566.2,  23754.26, 285,      37945.09, 79550.55)
    ,@("FL2792", 18340.02, 25513.84, 0,        37287.45, 81141.31)
    ,@("FL2900", 4049.5,   3892.18,  0,        12455.74, 20397.42)
    ,@("FL2948", 7695.15,  7055.77,  1180,     15603.58, 31534.5)
    ,@("FL2977", 36891.16, 55306.48, 0,        43520.25, 135717.89)
    ,@("FL3066", 3972.01,  5157.38,  0,        11046.75, 20176.14)
    ,@("FL3108", 4235.35,  6067.5,   0,        13255.03, 23557.88)
    ,@("FL3271", 39168.12, 59403.25, 0,        41576.81, 140148.18)
    ,@("FL3335", 4056.78,  4793.33,  11228.59, 12992.39, 33071.09)
)

$rowCount = $rows.Count
$colCount = 6

# Build a true 2D COM-friendly array (jagged arrays aren't accepted by Range.Value).
$arr = New-Object "object[,]" $rowCount, $colCount
for ($r = 0; $r -lt $rowCount; $r++) {
    for ($c = 0; $c -lt $colCount; $c++) {
        $arr[$r, $c] = $rows[$r][$c]
    }
}

$startRow = 2
$endRow = $startRow + $rowCount - 1
$target = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, $colCount))
$target.Value = $arr
